$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.840.06'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '3.129.43'
$ws.Range('E3').Value = '  +1.25%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '533.44'
$ws.Range('E5').Value = '  +1.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.93'
$ws.Range('E6').Value = '  +1.90%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '3.127.51'
$ws.Range('E8').Value = '  +1.19%  '
$ws.Range('E9').Value = '  +6.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.33'
$ws.Range('E10').Value = '  +0.34%  '
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.412'
$ws.Range('E12').Value = '  +4.86%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '3.666.16'
$ws.Range('E13').Value = '  +1.05%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.137'
$ws.Range('E14').Value = '  +1.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.93'
$ws.Range('E15').Value = '  +3.04%  '
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').Value = '57.936.28'
$ws.Range('E17').Value = '  +1.06%  '
$ws.Range('D18').Value = '3.126.71'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.09'
$ws.Range('E19').Value = '  +2.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.71'
$ws.Range('E20').Value = '  +3.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.10'
$ws.Range('E21').Value = '  +3.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '368.05'
$ws.Range('E22').Value = '  +6.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.68'
$ws.Range('E24').Value = '  -1.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '69.26'
$ws.Range('E25').Value = '  +2.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.507'
$ws.Range('E26').Value = '  +1.84%  '
$ws.Range('E27').Value = '  +1.12%  '
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('D29').Value = '0.0₃0866'
$ws.Range('E29').Value = '  -2.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.31'
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.88'
$ws.Range('E31').Value = '  +0.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.11'
$ws.Range('E32').Value = '  +1.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.48'
$ws.Range('E33').Value = '  +3.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.16'
$ws.Range('E34').Value = '  +5.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.18'
$ws.Range('E35').Value = '  +3.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.50'
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('E37').Value = '  +1.16%  '
$ws.Range('E38').Value = '  +5.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '25.51'
$ws.Range('E39').Value = '  -1.09%  '
$ws.Range('E40').Value = '  +5.20%  '
$ws.Range('E41').Value = '  +2.55%  '
$ws.Range('D42').Value = '2.530.64'
$ws.Range('E42').Value = '  +7.07%  '
$ws.Range('E43').Value = '  +0.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.700'
$ws.Range('E44').Value = '  +0.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '37.74'
$ws.Range('E45').Value = '  +3.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0270'
$ws.Range('E46').Value = '  +1.65%  '
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.979'
$ws.Range('E48').Value = '  +1.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.14'
$ws.Range('E49').Value = '  +3.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.79'
$ws.Range('E50').Value = '  +0.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.745'
$ws.Range('E51').Value = '  -0.77%  '
